$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: rename product and update price/quantity/total
$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1

# Row 4: new product "b" replaces the old "Sub Total" label position
$ws.Range("A4").Value = "b"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 4

# Row 5: now holds "Sub Total"
$ws.Range("A5").Value = "Sub Total"
$ws.Range("D5").Value = 5

# Row 6: now holds "HST"
$ws.Range("A6").Value = "HST"
$ws.Range("D6").Value = 0.65

# Row 7: new "Total" row
$ws.Range("A7").Value = "Total"
$ws.Range("D7").Value = 5.65
